# Upload Leave Card 12/27/2023 4:01 PM
# Adds a new "FL(1-0-0)"/"FL(3-0-0)" leave entry block for 2024 to the
# leave card table (Table1) on Sheet1, inserting 5 rows before the old
# row 89 and shifting all subsequent data down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Insert 5 new rows starting at row 89 (pushes old rows 89-116 down
#    to 94-121) and grow the Table1 listobject to match.
# ---------------------------------------------------------------------
$ws.Range("A89:A93").EntireRow.Insert()

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K121"))

# ---------------------------------------------------------------------
# Helper donor cells (untouched by the insert, i.e. above row 89) that
# already carry each cell style we need to stamp onto the new/changed
# rows, so every cell ends up with the exact same style index as the
# target workbook.
# ---------------------------------------------------------------------
$styleDonors = @{
    24 = "A9"
    41 = "A10"
    12 = "B57"
    21 = "B10"
    14 = "C9"
    40 = "D10"
    9  = "E58"
    39 = "K19"
}

function Set-Style([string]$cellRef, [int]$styleKey) {
    $donor = $styleDonors[$styleKey]
    $ws.Range($donor).Copy()
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats)
}

$IF_FORMULA = 'IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ---------------------------------------------------------------------
# 2. Row 86 : SL(1-0-0) full entry (date already present) - only the
#    REMARKS (K) cell style/value changes, rest of the row is filled in.
# ---------------------------------------------------------------------
Set-Style "K86" 39
$ws.Range("B86").Value = "SL(1-0-0)"
$ws.Range("C86").Value = 1.25
$ws.Range("G86").Formula = "=" + $IF_FORMULA
$ws.Range("H86").Value = 1
$ws.Range("K86").Value = 45201

# ---------------------------------------------------------------------
# 3. Row 87 : brand-new "split" sub-row (no date) for the same
#    SL(1-0-0) biweekly period.
# ---------------------------------------------------------------------
Set-Style "A87" 24
Set-Style "B87" 21
Set-Style "C87" 14
Set-Style "D87" 40
Set-Style "E87" 9
Set-Style "F87" 21
Set-Style "G87" 14
Set-Style "H87" 40
Set-Style "I87" 9
Set-Style "J87" 12
Set-Style "K87" 39

$ws.Range("A87").ClearContents()
$ws.Range("B87").Value = "SL(1-0-0)"
$ws.Range("C87").ClearContents()
$ws.Range("D87").ClearContents()
$ws.Range("G87").Formula = "=" + $IF_FORMULA
$ws.Range("H87").Value = 1
$ws.Range("K87").Value = 45230

# ---------------------------------------------------------------------
# 4. Row 88 : SL(1-0-0) full entry, analogous to row 86 (date moved up
#    from the old row 87 position).
# ---------------------------------------------------------------------
Set-Style "K88" 39
$ws.Range("A88").Value = 45231
$ws.Range("B88").Value = "SL(1-0-0)"
$ws.Range("C88").Value = 1.25
$ws.Range("G88").Formula = "=" + $IF_FORMULA
$ws.Range("H88").Value = 1
$ws.Range("K88").Value = 45238

# ---------------------------------------------------------------------
# 5. Row 89 : brand-new "split" sub-row, mirrors row 87.
# ---------------------------------------------------------------------
Set-Style "A89" 24
Set-Style "B89" 21
Set-Style "C89" 14
Set-Style "D89" 40
Set-Style "E89" 9
Set-Style "F89" 21
Set-Style "G89" 14
Set-Style "H89" 40
Set-Style "I89" 9
Set-Style "J89" 12
Set-Style "K89" 39

$ws.Range("A89").ClearContents()
$ws.Range("B89").Value = "SL(1-0-0)"
$ws.Range("C89").ClearContents()
$ws.Range("D89").ClearContents()
$ws.Range("G89").Formula = "=" + $IF_FORMULA
$ws.Range("H89").Value = 1
$ws.Range("K89").Value = 45247

# ---------------------------------------------------------------------
# 6. Row 90 : SL(1-0-0) entry without the EARNED amount (date moved up
#    from the old row 88 position); only K's style changes, rest keep
#    their existing formatting.
# ---------------------------------------------------------------------
Set-Style "K90" 39
$ws.Range("A90").Value = 45261
$ws.Range("B90").Value = "SL(1-0-0)"
$ws.Range("C90").ClearContents()
$ws.Range("G90").Formula = "=" + $IF_FORMULA
$ws.Range("H90").Value = 1
$ws.Range("K90").Value = 45271

# ---------------------------------------------------------------------
# 7. Row 91 : brand-new FL(1-0-0) entry (1 day).
# ---------------------------------------------------------------------
Set-Style "A91" 24
Set-Style "B91" 12
Set-Style "C91" 14
Set-Style "D91" 12
Set-Style "E91" 9
Set-Style "F91" 12
Set-Style "G91" 14
Set-Style "H91" 12
Set-Style "I91" 9
Set-Style "J91" 12
Set-Style "K91" 39

$ws.Range("A91").ClearContents()
$ws.Range("B91").Value = "FL(1-0-0)"
$ws.Range("C91").ClearContents()
$ws.Range("D91").Value = 1
$ws.Range("G91").ClearContents()
$ws.Range("H91").ClearContents()
$ws.Range("K91").Value = 45278

# ---------------------------------------------------------------------
# 8. Row 92 : brand-new FL(3-0-0) entry (3 days), REMARKS references
#    the 12/27-29/2023 leave dates.
# ---------------------------------------------------------------------
Set-Style "A92" 24
Set-Style "B92" 12
Set-Style "C92" 14
Set-Style "D92" 12
Set-Style "E92" 9
Set-Style "F92" 12
Set-Style "G92" 14
Set-Style "H92" 12
Set-Style "I92" 9
Set-Style "J92" 12
Set-Style "K92" 39

$ws.Range("A92").ClearContents()
$ws.Range("B92").Value = "FL(3-0-0)"
$ws.Range("C92").ClearContents()
$ws.Range("D92").Value = 3
$ws.Range("G92").ClearContents()
$ws.Range("H92").ClearContents()
$ws.Range("K92").Value = "12/27-29/2023"

# ---------------------------------------------------------------------
# 9. Row 93 : the "2024" year-label / totals row, same visual pattern
#    as the existing 2019/2021/2022/2023 label rows.
# ---------------------------------------------------------------------
Set-Style "A93" 41
Set-Style "B93" 12
Set-Style "C93" 14
Set-Style "D93" 12
Set-Style "E93" 9
Set-Style "F93" 12
Set-Style "G93" 14
Set-Style "H93" 12
Set-Style "I93" 9
Set-Style "J93" 12
Set-Style "K93" 39

$ws.Range("B93").ClearContents()
$ws.Range("C93").ClearContents()
$ws.Range("D93").ClearContents()
$ws.Range("G93").ClearContents()
$ws.Range("H93").ClearContents()
$ws.Range("K93").ClearContents()
$ws.Range("A93").Value = "2024"

# ---------------------------------------------------------------------
# 10. Recalculate so BALANCE (E9/I9) and the EARNED mirror column (G)
#     pick up the new rows, then restore the active selection.
# ---------------------------------------------------------------------
$excel.Calculate()

$ws.Range("E88").Select()
